# "Added clear map button" — the accompanying data change adds a new
# "Accident Type" column to the traffic-data sheet, inserted just before
# the existing "Summary" column (so Summary moves from column D to E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at D; everything at/after D (i.e. the old
# "Summary" column) shifts one column to the right.
$ws.Range("D1").EntireColumn.Insert()

# Populate the new "Accident Type" column.
$ws.Range("D1").Value = "Accident Type"
$ws.Range("D2").Value = "Animal"
$ws.Range("D3").Value = "Weather"

# Size the new column to fit its contents (matches the authored width).
$ws.Columns.Item(4).ColumnWidth = 11.46

# Leave the selection on D2, matching the saved workbook state.
$ws.Range("D2").Select()
